$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the most recent meeting row),
# pushing all existing meeting rows down by one.
$ws.Rows("2:2").Insert()

# Fill in the new row with the data for the newly added meeting (第641回).
$ws.Range("A2").Value = "第641回"
$ws.Range("B2").Value = "2026年1月14日（令和8年1月14日）"
$ws.Range("C2").Value = "１先進医療会議及び患者申出療養評価会議からの報告について`n２医療法等改正に伴う療養担当規則等の所要の見直しについて（諮問）`n３医療法等改正を踏まえた対応について（その２）`n４入院について（その９）`n５賃上げについて（その２）`n６物価対応について（その２）`n７これまでの議論の整理（案）について`n８令和８年度診療報酬改定について（諮問）`n９再生医療等製品の医療保険上の取扱いについて`n10費用対効果評価専門部会・薬価専門部会・保険医療材料専門部会 合同部会からの報告について`n11その他`n"
$ws.Range("D2").Value = "－"
$ws.Range("E2").Value = "資料`n"
$ws.Range("F2").Value = "－"
